$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# Delete the 36 rows of "New users of ..." drug cohort entries
# (rows 297 through 332 in the original sheet, cohort ids 1035-1070)
$ws.Range("A297:C332").EntireRow.Delete()
